{"js": "// Map of original text -> new text, per the diff.\nconst replacements = [\n  [\"2024-01-16 Tuesday\", \"2024-01-17 Wednesday\"],\n  [\"592\u00f78=74, 0\", \"731\u00f74=182, 3\"],\n  [\"843\u00f74=210, 3\", \"945\u00f72=472, 1\"],\n  [\"256\u00f78=32, 0\", \"702\u00f77=100, 2\"],\n  [\"803\u00f79=89, 2\", \"573\u00f79=63, 6\"],\n  [\"489\u00f74=122, 1\", \"904\u00f77=129, 1\"],\n  [\"974\u00f76=162, 2\", \"318\u00f74=79, 2\"],\n  [\"184\u00f72=92, 0\", \"774\u00f77=110, 4\"],\n  [\"567\u00f73=189, 0\", \"447\u00f79=49, 6\"],\n  [\"853\u00f73=284, 1\", \"653\u00f74=163, 1\"],\n  [\"971\u00f78=121, 3\", \"447\u00f77=63, 6\"],\n  [\"593\u00f79=65, 8\", \"822\u00f73=274, 0\"],\n  [\"576\u00f74=144, 0\", \"303\u00f74=75, 3\"],\n  [\"441\u00f78=55, 1\", \"835\u00f72=417, 1\"],\n  [\"994\u00f72=497, 0\", \"309\u00f77=44, 1\"],\n  [\"566\u00f78=70, 6\", \"162\u00f74=40, 2\"],\n  [\"233\u00f74=58, 1\", \"240\u00f73=80, 0\"],\n  [\"743\u00f78=92, 7\", \"527\u00f78=65, 7\"],\n  [\"266\u00f72=133, 0\", \"643\u00f75=128, 3\"],\n  [\"226\u00f72=113, 0\", \"350\u00f75=70, 0\"],\n  [\"305\u00f78=38, 1\", \"420\u00f77=60, 0\"],\n  [\"944\u00f75=188, 4\", \"530\u00f73=176, 2\"],\n  [\"296\u00f74=74, 0\", \"228\u00f74=57, 0\"],\n  [\"951\u00f73=317, 0\", \"929\u00f72=464, 1\"],\n  [\"361\u00f75=72, 1\", \"820\u00f79=91, 1\"],\n  [\"347\u00f72=173, 1\", \"136\u00f72=68, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five division problems/answers in the table,\n# per the commit's diff. Every \"before\" string in this document is unique, so a\n# plain Find/Replace (no wildcards) on each exact string is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2024-01-16 Tuesday\"; New = \"2024-01-17 Wednesday\"},\n    @{Old = \"592\u00f78=74, 0\";        New = \"731\u00f74=182, 3\"},\n    @{Old = \"843\u00f74=210, 3\";       New = \"945\u00f72=472, 1\"},\n    @{Old = \"256\u00f78=32, 0\";        New = \"702\u00f77=100, 2\"},\n    @{Old = \"803\u00f79=89, 2\";        New = \"573\u00f79=63, 6\"},\n    @{Old = \"489\u00f74=122, 1\";       New = \"904\u00f77=129, 1\"},\n    @{Old = \"974\u00f76=162, 2\";       New = \"318\u00f74=79, 2\"},\n    @{Old = \"184\u00f72=92, 0\";        New = \"774\u00f77=110, 4\"},\n    @{Old = \"567\u00f73=189, 0\";       New = \"447\u00f79=49, 6\"},\n    @{Old = \"853\u00f73=284, 1\";       New = \"653\u00f74=163, 1\"},\n    @{Old = \"971\u00f78=121, 3\";       New = \"447\u00f77=63, 6\"},\n    @{Old = \"593\u00f79=65, 8\";        New = \"822\u00f73=274, 0\"},\n    @{Old = \"576\u00f74=144, 0\";       New = \"303\u00f74=75, 3\"},\n    @{Old = \"441\u00f78=55, 1\";        New = \"835\u00f72=417, 1\"},\n    @{Old = \"994\u00f72=497, 0\";       New = \"309\u00f77=44, 1\"},\n    @{Old = \"566\u00f78=70, 6\";        New = \"162\u00f74=40, 2\"},\n    @{Old = \"233\u00f74=58, 1\";        New = \"240\u00f73=80, 0\"},\n    @{Old = \"743\u00f78=92, 7\";        New = \"527\u00f78=65, 7\"},\n    @{Old = \"266\u00f72=133, 0\";       New = \"643\u00f75=128, 3\"},\n    @{Old = \"226\u00f72=113, 0\";       New = \"350\u00f75=70, 0\"},\n    @{Old = \"305\u00f78=38, 1\";        New = \"420\u00f77=60, 0\"},\n    @{Old = \"944\u00f75=188, 4\";       New = \"530\u00f73=176, 2\"},\n    @{Old = \"296\u00f74=74, 0\";        New = \"228\u00f74=57, 0\"},\n    @{Old = \"951\u00f73=317, 0\";       New = \"929\u00f72=464, 1\"},\n    @{Old = \"361\u00f75=72, 1\";        New = \"820\u00f79=91, 1\"},\n    @{Old = \"347\u00f72=173, 1\";       New = \"136\u00f72=68, 0\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
